$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell reference -> new value map (columns D/E need text format to preserve exact string content)
$changes = [ordered]@{
    'D2' = '43.967.89'
    'E2' = '  +0.10%  '
    'D3' = '2.249.28'
    'E3' = '  -1.70%  '
    'E4' = '  +0.26%  '
    'D5' = '232.05'
    'E5' = '  -0.04%  '
    'D6' = '0.642'
    'E6' = '  +3.21%  '
    'D7' = '63.33'
    'E7' = '  -0.53%  '
    'E8' = '  +0.08%  '
    'E9' = '  +5.08%  '
    'D10' = '0.0980'
    'E10' = '  +3.15%  '
    'D11' = '57.34'
    'E11' = '  -0.61%  '
    'D12' = '26.31'
    'E12' = '  -0.26%  '
    'E13' = '  +1.42%  '
    'D14' = '2.582.94'
    'E14' = '  -1.66%  '
    'D15' = '15.51'
    'E15' = '  -2.17%  '
    'D16' = '6.12'
    'E16' = '  +2.53%  '
    'D17' = '0.830'
    'E17' = '  +1.62%  '
    'D18' = '2.247.01'
    'E18' = '  -1.55%  '
    'D19' = '43.832.76'
    'E19' = '  +0.12%  '
    'D20' = '0.0₃0986'
    'E20' = '  +3.33%  '
    'D21' = '72.68'
    'E21' = '  -0.78%  '
    'E22' = '  -2.08%  '
    'D23' = '248.07'
    'E23' = '  -1.69%  '
    'D25' = '2.41'
    'E25' = '  -7.45%  '
    'B26' = 'Toncoin'
    'C26' = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    'D26' = '2.30'
    'E26' = '  +0.86%  '
    'B27' = 'WEMIXToken'
    'C27' = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    'D27' = '3.34'
    'E27' = '  +21.04%  '
    'E28' = '  -0.75%  '
    'E29' = '  +0.21%  '
    'D30' = '21.00'
    'E30' = '  +2.02%  '
    'E31' = '  -0.73%  '
    'E32' = '  -2.60%  '
    'E33' = '  +2.17%  '
    'D34' = '0.0686'
    'E34' = '  -2.46%  '
    'E35' = '  +0.92%  '
    'D36' = '4.95'
    'E36' = '  -3.70%  '
    'D37' = '3.65'
    'E37' = '  -1.97%  '
    'D38' = '6.41'
    'E38' = '  -2.92%  '
    'E39' = '  -4.02%  '
    'D40' = '0.0253'
    'E40' = '  +1.18%  '
    'E41' = '  +0.02%  '
    'B42' = 'FraxShare'
    'C42' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D42' = '8.59'
    'E42' = '  -0.09%  '
    'B43' = 'TerraClassic'
    'C43' = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
    'D43' = '0.000224'
    'E43' = '  +1.35%  '
    'D44' = '17.08'
    'E44' = '  -0.21%  '
    'D45' = '97.39'
    'E45' = '  -1.17%  '
    'E46' = '  -2.63%  '
    'D47' = '0.0943'
    'E47' = '  -2.82%  '
    'D48' = '4.31'
    'E48' = '  -7.19%  '
    'D49' = '1.439.12'
    'E49' = '  -3.30%  '
    'D50' = '2.28'
    'E50' = '  -2.35%  '
    'E51' = '  +1.15%  '
}

foreach ($cellRef in $changes.Keys) {
    $col = ($cellRef -replace '[0-9]+$', '')
    $range = $ws.Range($cellRef)
    if ($col -eq 'D' -or $col -eq 'E') {
        $range.NumberFormat = "@"
    }
    $range.Value = $changes[$cellRef]
}
